$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (company id "4" -> "3"; refreshed financial metrics)
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "3"
$ws.Range("D2").Value = -0.10135
$ws.Range("E2").Value = 0.05769999999999999
$ws.Range("G2").Value = 0.3621467700702485
$ws.Range("H2").Value = 0.3621467700702485
$ws.Range("I2").Value = 0.3858389079028066
$ws.Range("J2").Value = 0.3276139867756158
$ws.Range("K2").Value = 726.05
$ws.Range("L2").Value = 0.07798771187350963
$ws.Range("M2").Value = 385.8
$ws.Range("N2").Value = 0.02185810925655233
$ws.Range("O2").Value = 0.5313683630603954
$ws.Range("P2").Value = 349
$ws.Range("Q2").Value = 0.01977314704649239
$ws.Range("R2").Value = 0.4806831485434888
$ws.Range("S2").Value = 36.80000000000001
$ws.Range("T2").Value = 0.09538621047174704
$ws.Range("U2").Value = 732.5
$ws.Range("V2").Value = 0.04150094616491597
$ws.Range("W2").Value = 0.1316033880674641
$ws.Range("X2").Value = 0.04110383326229912
$ws.Range("Y2").Value = 0.09049955480516501
$ws.Range("Z2").Value = 1.473279498234723
$ws.Range("AA2").Value = 0.255446318824245
$ws.Range("AB2").Value = 0.04007305844255434
$ws.Range("AC2").Value = 0.2152153955414566
$ws.Range("AD2").Value = 2229.408
$ws.Range("AE2").Value = 310.5346760322558
$ws.Range("AF2").Value = 2539.942676032256
$ws.Range("AG2").Value = 1807.442676032256
$ws.Range("AH2").Value = 0.1258011256674984
$ws.Range("AI2").Value = 0.2623706425252567
$ws.Range("AJ2").Value = 0.09289114339933104
$ws.Range("AK2").Value = 0.2019885626116808
$ws.Range("AL2").Value = 85.02499999999999
$ws.Range("AM2").Value = 85.02499999999999
$ws.Range("AN2").Value = 0.5871699584131139
$ws.Range("AO2").Value = 42.49679506027639
$ws.Range("AP2").Value = 0.4760349119227826
$ws.Range("AQ2").Value = 42.49679506027639

# Row 3 (Intact Financial Corporation - refreshed financial metrics)
$ws.Range("D3").Value = 0.0863
$ws.Range("E3").Value = 0.05769999999999999
$ws.Range("G3").Value = 0.3684158719048035
$ws.Range("H3").Value = 0.3684158719048035
$ws.Range("I3").Value = 0.3902236705741479
$ws.Range("J3").Value = 0.3167298322852374
$ws.Range("K3").Value = 708.5
$ws.Range("L3").Value = 0.07748928165193805
$ws.Range("M3").Value = 385.8
$ws.Range("N3").Value = 0.02280546196134066
$ws.Range("O3").Value = 0.5445306986591391
$ws.Range("P3").Value = 349
$ws.Range("Q3").Value = 0.0206301353667908
$ws.Range("R3").Value = 0.492589978828511
$ws.Range("S3").Value = 36.80000000000001
$ws.Range("T3").Value = 0.09538621047174704
$ws.Range("U3").Value = 628.2
$ws.Range("V3").Value = 0.03713424366022345
$ws.Range("W3").Value = 0.1316033880674641
$ws.Range("X3").Value = 0.04357350439231643
$ws.Range("Y3").Value = 0.0880298836751477
$ws.Range("Z3").Value = 1.466611456760101
$ws.Range("AA3").Value = 0.4645196007272345
$ws.Range("AB3").Value = 0.03976553175257089
$ws.Range("AC3").Value = 0.4247540689746636
$ws.Range("AD3").Value = 2200.6
$ws.Range("AE3").Value = 310.5346760322558
$ws.Range("AF3").Value = 2511.134676032256
$ws.Range("AG3").Value = 1882.934676032256
$ws.Range("AH3").Value = 0.1292524844976572
$ws.Range("AI3").Value = 0.2663617544060958
$ws.Range("AJ3").Value = 0.1001564478004692
$ws.Range("AK3").Value = 0.2139860279619797
$ws.Range("AL3").Value = 84.1
$ws.Range("AM3").Value = 84.1
$ws.Range("AN3").Value = 0.5835277895630038
$ws.Range("AO3").Value = 42.67657550535078
$ws.Range("AP3").Value = 0.4992932424777937
$ws.Range("AQ3").Value = 42.67657550535078

# Row 4 (company swapped: Trisura Group Ltd. -> EFH Holdings Inc.)
$ws.Range("B4").Value = "EFH Holdings Inc. (TSXV:EFH)"
$ws.Range("D4").Value = -0.289
$ws.Range("G4").Value = -0.01664473684210526
$ws.Range("H4").Value = -0.01664473684210526
$ws.Range("I4").Value = 0.124671052631579
$ws.Range("J4").Value = 0.1130602045980336
$ws.Range("K4").Value = -1.75
$ws.Range("L4").Value = -0.05756578947368422
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("U4").Value = 10.6
$ws.Range("V4").Value = 0.6883116883116883
$ws.Range("W4").Value = -0.0258493353028065
$ws.Range("X4").Value = 0.04065628174844418
$ws.Range("Y4").Value = -0.06650561705125069
$ws.Range("Z4").Value = 2.259383128948344
$ws.Range("AA4").Value = 0.255446318824245
$ws.Range("AB4").Value = 0.04023092328278841
$ws.Range("AC4").Value = 0.2152153955414566
$ws.Range("AD4").Value = 0.308
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 0.308
$ws.Range("AG4").Value = -10.292
$ws.Range("AH4").Value = 0.0196078431372549
$ws.Range("AI4").Value = 0.01749204906860518
$ws.Range("AJ4").Value = -2.014878621769773
$ws.Range("AK4").Value = -1.468607305936073
$ws.Range("AL4").Value = 0
$ws.Range("AM4").Value = 0
$ws.Range("AN4").Value = 0.07758186397984886
$ws.Range("AP4").Value = -2.592443324937027
$ws.Range("AO4").ClearContents()
$ws.Range("AQ4").ClearContents()

# Row 5 (company swapped: Echelon Financial Holdings Inc. -> Trisura Group Ltd.)
$ws.Range("B5").Value = "Trisura Group Ltd. (TSX:TSU)"
$ws.Range("G5").Value = 0.02584434654919237
$ws.Range("H5").Value = 0.02584434654919237
$ws.Range("I5").Value = 0.1497797356828194
$ws.Range("J5").Value = 0.1241307594864911
$ws.Range("K5").Value = 19.3
$ws.Range("L5").Value = 0.1417033773861968
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 0
$ws.Range("U5").Value = 93.7
$ws.Range("V5").Value = 0.1305377542490945
$ws.Range("W5").Value = 0.1364922206506365
$ws.Range("X5").Value = 0.04110383326229912
$ws.Range("Y5").Value = 0.09538838738833737
$ws.Range("Z5").Value = 1.907295896933202
$ws.Range("AA5").Value = 0.2367540882517866
$ws.Range("AB5").Value = 0.04007305844255434
$ws.Range("AC5").Value = 0.1966810298092322
$ws.Range("AD5").Value = 28.5
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 28.5
$ws.Range("AG5").Value = -65.2
$ws.Range("AH5").Value = 0.03818839608736433
$ws.Range("AI5").Value = 0.1209677419354839
$ws.Range("AJ5").Value = -0.09990806006742263
$ws.Range("AK5").Value = -0.459478505990134
$ws.Range("AL5").Value = 0.925
$ws.Range("AM5").Value = 0.925
$ws.Range("AN5").Value = 1.313364055299539
$ws.Range("AO5").Value = 22.05405405405405
$ws.Range("AP5").Value = -3.004608294930876
$ws.Range("AQ5").Value = 22.05405405405405
$ws.Range("D5").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("T5").ClearContents()

# Row 6 (Kingsway Financial Services Inc.) removed entirely
$ws.Rows.Item(6).Delete()
